$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out old rows 4, 7, 8 (B column values, and shift rows up to 2 and 3)
$ws.Range("A4:B8").ClearContents()

# Write new data: row2 = Super Strike, row3 = LVE Mission (only column A, no B values)
$ws.Range("A2").Value = "Super Strike"
$ws.Range("A3").Value = "LVE Mission"

# Update selection to match the target state
$ws.Range("B3").Select()
